# Insert a new weekly price record for "Berenjena" (Macroferia Regional de Talca)
# as row 192 of Sheet1, pushing the existing rows 192:212 down to 193:213.
#
# This mirrors the target diff: dimension grows from A1:R212 to A1:R213, and the
# new row 192 carries the newest reading (Fecha serial 45180 = 2023-09-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 192..212 down to 193..213 (inherits formatting, incl. the date
# style on column D, from the row above, same as a manual Excel row-insert).
$ws.Rows("192:192").Insert()

# Populate the newly inserted row 192 with the new record's data.
$ws.Cells.Item(192, 1).Value  = 5
$ws.Cells.Item(192, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(192, 3).Value  = "Maule"
$ws.Cells.Item(192, 4).Value  = 45180
$ws.Cells.Item(192, 5).Value  = 7
$ws.Cells.Item(192, 6).Value  = 100112001
$ws.Cells.Item(192, 7).Value  = "Berenjena"
$ws.Cells.Item(192, 8).Value  = "Sin especificar"
$ws.Cells.Item(192, 9).Value  = "Primera"
$ws.Cells.Item(192, 10).Value = 200
$ws.Cells.Item(192, 11).Value = 10000
$ws.Cells.Item(192, 12).Value = 10000
$ws.Cells.Item(192, 13).Value = 10000
$ws.Cells.Item(192, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(192, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(192, 16).Value = 200
$ws.Cells.Item(192, 17).Value = 50
$ws.Cells.Item(192, 18).Value = "Hortaliza"
